# Update the "Förändrad" (Changed) date column (C) for rows 2 through 18
# from serial date 45192 (2023-09-23) to 45202 (2023-10-03).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("C2:C18")
$rng.Value2 = 45202
